$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending / target cluster labels, in order, cross-joined to rebuild rows 2-10
# (3 sending clusters x 3 target clusters), per "Natmi following Dr Hou advice".
$clusters = @("ECs", "FAPs", "sCs")

# Per-(sending,target) pair metric values, columns E..T (16 numeric columns).
# Row order follows sending cluster outer loop, target cluster inner loop
# (ECs->ECs, ECs->FAPs, ECs->sCs, FAPs->ECs, FAPs->FAPs, FAPs->sCs, sCs->ECs, sCs->FAPs, sCs->sCs).
$data = @(
    @(3, 1, 0.6873773333333334, 2.062132, 0.02660947569874856, 0.02660947569874856, 2, 0.6666666666666666, 0.09433999999999999, 0.28302, 0.05191071108246543, 0.05191071108246543, 0.06484717762666667, 0.58362459864, 0.001381316805053621, 0.001381316805053621),
    @(3, 1, 0.6873773333333334, 2.062132, 0.02660947569874856, 0.02660947569874856, 3, 1, 0.9431116666666667, 2.829335, 0.5189484550226392, 0.5189484550226391, 0.648273582468889, 5.83446224222, 0.01380894630282803, 0.01380894630282802),
    @(3, 1, 0.6873773333333334, 2.062132, 0.02660947569874856, 0.02660947569874856, 3, 1, 0.7798996666666667, 2.339699, 0.4291408338948954, 0.4291408338948954, 0.5360853531408889, 4.824768178268, 0.01141921259086691, 0.01141921259086691),
    @(3, 1, 5.913984666666667, 17.741954, 0.2289398029860915, 0.2289398029860915, 2, 0.6666666666666666, 0.09433999999999999, 0.28302, 0.05191071108246543, 0.05191071108246543, 0.5579253134533333, 5.02132782108, 0.01188442796808755, 0.01188442796808755),
    @(3, 1, 5.913984666666667, 17.741954, 0.2289398029860915, 0.2289398029860915, 3, 1, 0.9431116666666667, 2.829335, 0.5189484550226392, 0.5189484550226391, 5.577547935621111, 50.19793142059, 0.1188079570528196, 0.1188079570528195),
    @(3, 1, 5.913984666666667, 17.741954, 0.2289398029860915, 0.2289398029860915, 3, 1, 0.7798996666666667, 2.339699, 0.4291408338948954, 0.4291408338948954, 4.612314670205111, 41.510832031846, 0.09824741796518435, 0.09824741796518435),
    @(3, 1, 19.230689, 57.692067, 0.7444507213151601, 0.7444507213151601, 2, 0.6666666666666666, 0.09433999999999999, 0.28302, 0.05191071108246543, 0.05191071108246543, 1.81422320026, 16.32800880234, 0.03864496630932426, 0.03864496630932426),
    @(3, 1, 19.230689, 57.692067, 0.7444507213151601, 0.7444507213151601, 3, 1, 0.9431116666666667, 2.829335, 0.5189484550226392, 0.5189484550226391, 18.13668715393834, 163.230184385445, 0.3863315516669917, 0.3863315516669916),
    @(3, 1, 19.230689, 57.692067, 0.7444507213151601, 0.7444507213151601, 3, 1, 0.7798996666666667, 2.339699, 0.4291408338948954, 0.4291408338948954, 14.99800794087033, 134.982071467833, 0.3194742033388442, 0.3194742033388442)
)

$row = 2
$idx = 0
foreach ($sending in $clusters) {
    foreach ($target in $clusters) {
        $ws.Cells.Item($row, 1).Value = $sending   # A: Sending cluster
        $ws.Cells.Item($row, 2).Value = "Fgf1"      # B: Ligand symbol
        $ws.Cells.Item($row, 3).Value = "Fgfr2"     # C: Receptor symbol
        $ws.Cells.Item($row, 4).Value = $target     # D: Target cluster

        $vals = $data[$idx]
        for ($c = 0; $c -lt $vals.Length; $c++) {
            $ws.Cells.Item($row, 5 + $c).Value = $vals[$c]
        }

        $row = $row + 1
        $idx = $idx + 1
    }
}
